$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (font/border/alignment) from H1 into I1 and J1
# xlPasteFormats = -4122
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Set the new header text
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in data rows 2-13: column I = 1, column J = same value as column H
for ($r = 2; $r -le 13; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
